$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had an unused, style-only column A (blank) with the real data
# living in columns B:D. Remove that spacer column so the data shifts left
# to A:C. ClearFormats() first so the deleted column's formatting doesn't
# leave a stray zero-width <col> entry behind.
$ws.Columns.Item(1).ClearFormats()
$ws.Columns.Item(1).Delete()

# Add Benjamin Best's note about the "Aug " header-detection heuristic on
# what is now cell A2 (the "Aug 6, day 1, Mon: ..." row header).
$commentText = "Benjamin Best:" + [char]10 + "Agenda looks for Times starting with `"Aug `" to determine if inserting header."
$comment = $ws.Range("A2").AddComment($commentText)

# Match the post-edit selection/active cell.
$ws.Range("B10").Select() | Out-Null
